$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.375.17"
$ws.Range("E2").Value = "  +4.64%  "
$ws.Range("D3").Value = "4.047.67"
$ws.Range("E3").Value = "  +4.14%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'521.88"
$ws.Range("E5").Value = "  -1.43%  "
$ws.Range("D6").Value = "'148.72"
$ws.Range("E6").Value = "  +3.21%  "
$ws.Range("D7").Value = "'0.720"
$ws.Range("E7").Value = "  +17.84%  "
$ws.Range("D8").Value = "4.039.15"
$ws.Range("E8").Value = "  +4.21%  "
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("D10").Value = "'0.781"
$ws.Range("E10").Value = "  +8.80%  "
$ws.Range("D11").Value = "'0.178"
$ws.Range("E11").Value = "  +4.19%  "
$ws.Range("D12").Value = "'0.0000333"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "'48.48"
$ws.Range("E13").Value = "  +15.70%  "
$ws.Range("D14").Value = "'11.20"
$ws.Range("E14").Value = "  +9.83%  "
$ws.Range("D15").Value = "4.694.67"
$ws.Range("E15").Value = "  +4.13%  "
$ws.Range("D16").Value = "4.045.81"
$ws.Range("E16").Value = "  +4.07%  "
$ws.Range("D17").Value = "'21.31"
$ws.Range("E17").Value = "  +4.99%  "
$ws.Range("D18").Value = "'14.34"
$ws.Range("E18").Value = "  +2.57%  "
$ws.Range("D19").Value = "'1.22"
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("D20").Value = "'0.133"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").Value = "72.327.00"
$ws.Range("E21").Value = "  +4.63%  "
$ws.Range("D22").Value = "'446.70"
$ws.Range("E22").Value = "  +5.55%  "
$ws.Range("D23").Value = "'105.22"
$ws.Range("E23").Value = "  +20.32%  "
$ws.Range("D24").Value = "'3.62"
$ws.Range("E24").Value = "  +7.01%  "
$ws.Range("D25").Value = "'15.19"
$ws.Range("E25").Value = "  +7.31%  "
$ws.Range("D26").Value = "'4.04"
$ws.Range("E26").Value = "  +0.78%  "
$ws.Range("D27").Value = "'11.55"
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("D28").Value = "'11.21"
$ws.Range("E28").Value = "  +6.12%  "
$ws.Range("D29").Value = "'38.15"
$ws.Range("E29").Value = "  +5.13%  "
$ws.Range("D30").Value = "'5.82"
$ws.Range("E30").Value = "  +2.56%  "
$ws.Range("D31").Value = "'3.28"
$ws.Range("E31").Value = "  +15.07%  "
$ws.Range("D32").Value = "'13.83"
$ws.Range("E32").Value = "  +4.92%  "
$ws.Range("D33").Value = "'0.132"
$ws.Range("E33").Value = "  +4.33%  "
$ws.Range("D34").Value = "'680.66"
$ws.Range("E34").Value = "  -2.22%  "
$ws.Range("D35").Value = "'6.74"
$ws.Range("E35").Value = "  +13.50%  "
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("D37").Value = "'42.88"
$ws.Range("E37").Value = "  +7.05%  "
$ws.Range("E38").Value = "  +2.76%  "
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("E40").Value = "  +3.90%  "
$ws.Range("D41").Value = "'3.50"
$ws.Range("E41").Value = "  +6.65%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D43").Value = "'0.0501"
$ws.Range("E43").Value = "  +4.27%  "
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").Value = "'3.21"
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("E46").Value = "  +12.51%  "
$ws.Range("E47").Value = "  -1.83%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").Value = "'9.80"
$ws.Range("E48").Value = "  +14.79%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").Value = "'3.44"
$ws.Range("E49").Value = "  +0.88%  "
$ws.Range("E50").Value = "  +3.14%  "
$ws.Range("D51").Value = "'0.000278"
$ws.Range("E51").Value = "  +3.76%  "

Write-Host "Updated cryptos list"
